# Progress update as of 04-Nov-2025:
#  - Remove the expired "Misaligned Pallet-ULD Recovery Procedure (SOPs)" row
#    (row 14), shifting all following rows up by one.
#  - Decrement "PERIOD TO EXPIRE" (col H) by 1 day and bump "LAST UPDATE"
#    (col I) from 03-Nov-2025 to 04-Nov-2025 for every remaining data row.
#  - Renumber the "SN" column (col A) sequentially for the rows that shifted.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Drop the now-expired SOP row; everything below shifts up one row.
$ws.Rows.Item(14).Delete()

# New SN values (col A) for the rows that moved up after the deletion.
$snFixes = @{
    14 = 12
    15 = 13
    16 = 14
}
foreach ($row in $snFixes.Keys) {
    $ws.Cells.Item($row, 1).Value = $snFixes[$row]
}

# New "PERIOD TO EXPIRE" (col H) values, one day less than before.
$periodFixes = @{
    3  = 598
    4  = 277
    5  = 680
    6  = 678
    7  = 374
    8  = 444
    9  = 621
    10 = 678
    11 = 678
    12 = -30
    13 = 137
    14 = 136
    15 = 151
    16 = 607
}
foreach ($row in $periodFixes.Keys) {
    $ws.Cells.Item($row, 8).Value = $periodFixes[$row]
}

# "LAST UPDATE" (col I) moves from 03-Nov-2025 to 04-Nov-2025 for every
# remaining data row (3..16). Leading apostrophe keeps it literal text
# (matching the source file) instead of Excel auto-converting it to a
# date serial.
for ($row = 3; $row -le 16; $row++) {
    $ws.Cells.Item($row, 9).Value = "'04-Nov-2025"
}
